# Variant 8, problem 23: switch the dynamic-programming formulas in A15:L24
# from MIN(...) to MAX(...) (keeping the shared-formula structure / cell
# references intact), then move the active selection to A15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A15:L24")
[void]$rng.Replace("MIN(", "MAX(", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)

$excel.Calculate()

# Move the active cell selection from S13 to A15
[void]$ws.Range("A15").Select()
